$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Active Signals" sheet: grows from 2 data rows to 5 data rows (A1:J3 -> A1:J6)
# ---------------------------------------------------------------------------
$wsActive = $wb.Worksheets.Item("Active Signals")

# Row 2 currently carries the BUY (green) style, row 3 the SELL (red) style.
# Use them as templates so the new rows inherit the correct fill/border
# formatting, then overwrite every cell's value afterwards.
$wsActive.Range("A3:J3").Copy($wsActive.Range("A4:J4"))
$wsActive.Range("A3:J3").Copy($wsActive.Range("A5:J5"))
$wsActive.Range("A2:J2").Copy($wsActive.Range("A6:J6"))

# Row 2: XAUCHF BUY
$wsActive.Cells.Item(2,1).Value = "2025-07-28 20:03"
$wsActive.Cells.Item(2,2).Value = "XAUCHF"
$wsActive.Cells.Item(2,3).Value = "BUY"
$wsActive.Cells.Item(2,4).Value = 2334.28355
$wsActive.Cells.Item(2,5).Value = 2334.27912
$wsActive.Cells.Item(2,6).Value = 2334.29131
$wsActive.Cells.Item(2,7).Value = 0.06
$wsActive.Cells.Item(2,8).Value = "84.0%"
$wsActive.Cells.Item(2,9).Value = 1.75
$wsActive.Cells.Item(2,10).Value = "Active"

# Row 3: NZDUSD BUY
$wsActive.Cells.Item(3,1).Value = "2025-07-28 19:51"
$wsActive.Cells.Item(3,2).Value = "NZDUSD"
$wsActive.Cells.Item(3,3).Value = "BUY"
$wsActive.Cells.Item(3,4).Value = 0.59221
$wsActive.Cells.Item(3,5).Value = 0.58799
$wsActive.Cells.Item(3,6).Value = 0.59656
$wsActive.Cells.Item(3,7).Value = 0.02
$wsActive.Cells.Item(3,8).Value = "76.0%"
$wsActive.Cells.Item(3,9).Value = 1.03
$wsActive.Cells.Item(3,10).Value = "Active"

# Row 4: XAUCHF SELL
$wsActive.Cells.Item(4,1).Value = "2025-07-28 19:43"
$wsActive.Cells.Item(4,2).Value = "XAUCHF"
$wsActive.Cells.Item(4,3).Value = "SELL"
$wsActive.Cells.Item(4,4).Value = 2336.548
$wsActive.Cells.Item(4,5).Value = 2336.55109
$wsActive.Cells.Item(4,6).Value = 2336.54131
$wsActive.Cells.Item(4,7).Value = 0.07000000000000001
$wsActive.Cells.Item(4,8).Value = "87.0%"
$wsActive.Cells.Item(4,9).Value = 2.16
$wsActive.Cells.Item(4,10).Value = "Active"

# Row 5: USDJPY SELL
$wsActive.Cells.Item(5,1).Value = "2025-07-28 20:08"
$wsActive.Cells.Item(5,2).Value = "USDJPY"
$wsActive.Cells.Item(5,3).Value = "SELL"
$wsActive.Cells.Item(5,4).Value = 149.07482
$wsActive.Cells.Item(5,5).Value = 149.36232
$wsActive.Cells.Item(5,6).Value = 148.34779
$wsActive.Cells.Item(5,7).Value = 0.04
$wsActive.Cells.Item(5,8).Value = "85.0%"
$wsActive.Cells.Item(5,9).Value = 2.53
$wsActive.Cells.Item(5,10).Value = "Active"

# Row 6: USDJPY BUY
$wsActive.Cells.Item(6,1).Value = "2025-07-28 20:24"
$wsActive.Cells.Item(6,2).Value = "USDJPY"
$wsActive.Cells.Item(6,3).Value = "BUY"
$wsActive.Cells.Item(6,4).Value = 149.10511
$wsActive.Cells.Item(6,5).Value = 148.847
$wsActive.Cells.Item(6,6).Value = 150.01508
$wsActive.Cells.Item(6,7).Value = 0.09
$wsActive.Cells.Item(6,8).Value = "81.0%"
$wsActive.Cells.Item(6,9).Value = 3.53
$wsActive.Cells.Item(6,10).Value = "Active"

# ---------------------------------------------------------------------------
# 2) "Summary Dashboard" sheet: refresh the rolled-up stats
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary Dashboard")
$wsSummary.Cells.Item(4,2).Value = 5
$wsSummary.Cells.Item(7,2).Value = "85.2%"
$wsSummary.Cells.Item(8,2).Value = "2.07"
$wsSummary.Cells.Item(9,2).Value = "2025-07-28 20:00:31"

# ---------------------------------------------------------------------------
# 3) "Signal History" sheet: the log is reshuffled/updated in place (still 15 rows)
# ---------------------------------------------------------------------------
$wsHistory = $wb.Worksheets.Item("Signal History")

function Set-HistoryRow {
    param($ws, $r, $ts, $sym, $sig, $entry, $sl, $tp, $lots, $conf, $rr, $status)
    $ws.Cells.Item($r,1).Value = $ts
    $ws.Cells.Item($r,2).Value = $sym
    $ws.Cells.Item($r,3).Value = $sig
    $ws.Cells.Item($r,4).Value = $entry
    $ws.Cells.Item($r,5).Value = $sl
    $ws.Cells.Item($r,6).Value = $tp
    $ws.Cells.Item($r,7).Value = $lots
    $ws.Cells.Item($r,8).Value = $conf
    $ws.Cells.Item($r,9).Value = $rr
    $ws.Cells.Item($r,10).Value = $status
}

Set-HistoryRow $wsHistory 2  "2025-07-28 20:10" "USDJPY" "BUY"  149.43539   149.15827   150.33336   0.04                 0.83                 3.24 "Filled"
Set-HistoryRow $wsHistory 3  "2025-07-28 20:10" "XAUGBP" "SELL" 2105.39076  2105.39337  2105.38149  0.05                 0.9                  3.55 "Filled"
Set-HistoryRow $wsHistory 4  "2025-07-28 20:03" "XAUCHF" "BUY"  2334.28355  2334.27912  2334.29131  0.06                 0.84                 1.75 "Active"
Set-HistoryRow $wsHistory 5  "2025-07-28 19:51" "NZDUSD" "BUY"  0.59221     0.58799     0.59656     0.02                 0.76                 1.03 "Active"
Set-HistoryRow $wsHistory 6  "2025-07-28 20:01" "XAUAUD" "BUY"  4068.58247  4068.57785  4068.5889   0.05                 0.91                 1.39 "Filled"
Set-HistoryRow $wsHistory 7  "2025-07-28 19:43" "XAUCHF" "SELL" 2336.548    2336.55109  2336.54131  0.07000000000000001  0.87                 2.16 "Active"
Set-HistoryRow $wsHistory 8  "2025-07-28 19:41" "EURUSD" "BUY"  1.10507     1.1005      1.10936     0.03                 0.83                 0.9399999999999999 "Pending"
Set-HistoryRow $wsHistory 9  "2025-07-28 20:19" "XAUUSD" "BUY"  2644.48224  2644.47753  2644.48859  0.05                 0.95                 1.35 "Filled"
Set-HistoryRow $wsHistory 10 "2025-07-28 19:44" "NZDUSD" "SELL" 0.58648     0.58863     0.58035     0.02                 0.9399999999999999  2.85 "Filled"
Set-HistoryRow $wsHistory 11 "2025-07-28 20:08" "USDJPY" "SELL" 149.07482   149.36232   148.34779   0.04                 0.85                 2.53 "Active"
Set-HistoryRow $wsHistory 12 "2025-07-28 20:25" "USDCAD" "SELL" 1.36369     1.36737     1.35429     0.09                 0.77                 2.56 "Filled"
Set-HistoryRow $wsHistory 13 "2025-07-28 20:24" "USDJPY" "BUY"  149.10511   148.847     150.01508   0.09                 0.8100000000000001  3.53 "Active"
Set-HistoryRow $wsHistory 14 "2025-07-28 19:51" "XAUEUR" "SELL" 2422.95788  2422.96252  2422.95307  0.07000000000000001  0.84                 1.04 "Pending"
Set-HistoryRow $wsHistory 15 "2025-07-28 19:57" "XAUUSD" "SELL" 2649.17888  2649.18361  2649.17293  0.06                 0.78                 1.26 "Filled"
Set-HistoryRow $wsHistory 16 "2025-07-28 20:21" "USDCHF" "BUY"  0.88436     0.87957     0.89366     0.07000000000000001  0.9                  1.94 "Filled"
